# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $ws.Range($col + $Row).Value = $Values[$col]
    }
}

# --- Rows 104 and 105 swap their match data (ids, odds, etc.) ---
Set-RowValues 104 @{
    B  = 7127370
    F  = "Macarthur FC"
    G  = "Wellington Phoenix"
    H  = 1
    I  = 2
    J  = "A"
    K  = 2.4
    L  = 3.75
    M  = 2.625
    N  = 2.375
    O  = 3.8
    P  = 2.75
    Q  = 0
    R  = 1.8
    S  = 2.05
    T  = 3
    U  = 1.9
    V  = 1.95
    W  = -1
    X  = -1
    Y  = 1.75
    Z  = -1
    AA = 1.05
    AB = 0
    AC = -0
}

Set-RowValues 105 @{
    B  = 7127374
    F  = "Central Coast Mariners"
    G  = "Western Sydney Wanderers"
    H  = 1
    I  = 0
    J  = "H"
    K  = 1.909
    L  = 3.75
    M  = 3.6
    N  = 2.15
    O  = 3.6
    P  = 3.25
    Q  = -0.25
    R  = 1.86
    S  = 2.04
    T  = 2.75
    U  = 1.975
    V  = 1.875
    W  = 1.15
    X  = -1
    Y  = -1
    Z  = 0.8600000000000001
    AA = -1
    AB = -1
    AC = 0.875
}

# --- Rows 112 and 113 swap their match data ---
Set-RowValues 112 @{
    B  = 7127376
    F  = "Newcastle Jets"
    G  = "Macarthur FC"
    H  = 2
    I  = 2
    J  = "D"
    K  = 1.95
    L  = 4
    M  = 3.4
    N  = 1.909
    O  = 4.2
    P  = 3.6
    Q  = -0.5
    R  = 1.89
    S  = 2.01
    T  = 3.5
    U  = 1.95
    V  = 1.9
    W  = -1
    X  = 3.2
    Y  = -1
    Z  = -1
    AA = 1.01
    AB = 0.95
    AC = -1
}

Set-RowValues 113 @{
    B  = 7127379
    F  = "Melbourne Victory"
    G  = "Central Coast Mariners"
    H  = 0
    I  = 1
    J  = "A"
    K  = 1.95
    L  = 3.6
    M  = 3.8
    N  = 1.909
    O  = 3.6
    P  = 4
    Q  = -0.5
    R  = 1.9
    S  = 1.95
    T  = 2.75
    U  = 1.925
    V  = 1.925
    W  = -1
    X  = -1
    Y  = 3
    Z  = -1
    AA = 0.95
    AB = -1
    AC = 0.925
}

# --- Rows 124 and 125 swap their match data ---
Set-RowValues 124 @{
    B  = 7127388
    F  = "Sydney FC"
    G  = "Brisbane Roar"
    H  = 1
    I  = 1
    J  = "D"
    K  = 1.5
    L  = 5
    M  = 5
    N  = 1.533
    O  = 5.25
    P  = 5
    Q  = -1
    R  = 1.8
    S  = 2.05
    T  = 3.5
    U  = 1.925
    V  = 1.925
    W  = -1
    X  = 4.25
    Y  = -1
    Z  = -1
    AA = 1.05
    AB = -1
    AC = 0.925
}

Set-RowValues 125 @{
    B  = 7128012
    F  = "Macarthur FC"
    G  = "Central Coast Mariners"
    H  = 0
    I  = 3
    J  = "A"
    K  = 2.4
    L  = 3.5
    M  = 2.75
    N  = 3.4
    O  = 3.75
    P  = 2.05
    Q  = 0.25
    R  = 2.025
    S  = 1.825
    T  = 3
    U  = 2.05
    V  = 1.8
    W  = -1
    X  = -1
    Y  = 1.05
    Z  = -1
    AA = 0.825
    AB = 0
    AC = -0
}

# --- Rows 139-144 refreshed with updated fixtures/odds (old row 139 removed, 140-145 shift up) ---
Set-RowValues 139 @{
    B = 7127399
    E = 45387.23958333334
    F = "Western Sydney Wanderers"
    G = "Brisbane Roar"
    K = 2.1
    L = 3.75
    M = 3.1
    N = 1.95
    O = 4
    P = 3.5
    Q = -0.5
    R = 1.98
    S = 1.92
    T = 3.25
    U = 2
    V = 1.85
}

Set-RowValues 140 @{
    B = 8034339
    E = 45388.0625
    F = "Western United FC"
    G = "Macarthur FC"
    K = 2.6
    L = 3.5
    M = 2.55
    N = 2.05
    O = 3.8
    P = 3.1
    Q = -0.25
    R = 1.85
    S = 2.05
    T = 3.25
    U = 1.825
    V = 2.025
}

Set-RowValues 141 @{
    B = 8005739
    E = 45388.14583333334
    F = "Central Coast Mariners"
    G = "Wellington Phoenix"
    K = 1.8
    L = 3.6
    M = 4.333
    N = 1.727
    O = 3.6
    P = 4.75
    Q = -0.75
    R = 1.99
    S = 1.91
    T = 2.75
    U = 2
    V = 1.85
}

Set-RowValues 142 @{
    B = 7126794
    E = 45388.23958333334
    F = "Melbourne Victory"
    G = "Melbourne City"
    K = 1.833
    L = 3.5
    M = 4.5
    N = 2.2
    O = 3.75
    P = 3
    Q = -0.25
    R = 1.99
    S = 1.91
    T = 2.75
    U = 1.8
    V = 2.05
}

Set-RowValues 143 @{
    B = 7127403
    E = 45389.08333333334
    F = "Newcastle Jets"
    G = "Sydney FC"
    K = 3.6
    L = 3.6
    M = 1.952
    N = 4
    O = 3.8
    P = 1.8
    Q = 0.5
    R = 2.05
    S = 1.85
    T = 3.25
    U = 2
    V = 1.85
}

Set-RowValues 144 @{
    B = 7127402
    E = 45389.16666666666
    F = "Perth Glory"
    G = "Adelaide United"
    K = 2.25
    L = 3.5
    M = 3
    N = 2.625
    O = 3.5
    P = 2.55
    Q = 0
    R = 2.06
    S = 1.84
    T = 3.5
    U = 2.025
    V = 1.825
}

# --- Old row 145 no longer exists; delete it so the sheet dimension becomes A1:AC144 ---
$ws.Rows.Item(145).Delete()
